$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the workers table (rows 16-27) so that records are grouped by
# worker (Periodo Mora descending: 1806, 1805, 1804) instead of grouped
# by period. The set of 12 records stays the same, only the ordering of
# rows changes.

$data = @(
    @("34942765",   "YOMAIRA ELENA BANQUEZ SOTO",      "1806", 38674, 966842),
    @("34942765",   "YOMAIRA ELENA BANQUEZ SOTO",      "1805", 38674, 966842),
    @("34942765",   "YOMAIRA ELENA BANQUEZ SOTO",      "1804", 38674, 966842),
    @("1001969862", "ISAAC FRANCISCO FRANCIS OSPINA",  "1806", 65548, 1638690),
    @("1001969862", "ISAAC FRANCISCO FRANCIS OSPINA",  "1805", 65548, 1638690),
    @("1001969862", "ISAAC FRANCISCO FRANCIS OSPINA",  "1804", 65548, 1638690),
    @("45493646",   "DUBYS OROZCO RODRIGUEZ",          "1806", 54053, 1351320),
    @("45493646",   "DUBYS OROZCO RODRIGUEZ",          "1805", 54053, 1351320),
    @("45493646",   "DUBYS OROZCO RODRIGUEZ",          "1804", 54053, 1351320),
    @("22801823",   "CLAUDIA MARIA GOMEZ ZURIQUEZ",    "1806", 54053, 1351320),
    @("22801823",   "CLAUDIA MARIA GOMEZ ZURIQUEZ",    "1805", 54053, 1351320),
    @("22801823",   "CLAUDIA MARIA GOMEZ ZURIQUEZ",    "1804", 54053, 1351320)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 3).Value = $rec[0]
    $ws.Cells.Item($row, 4).Value = $rec[1]
    $ws.Cells.Item($row, 5).Value = $rec[2]
    $ws.Cells.Item($row, 6).Value = $rec[3]
    $ws.Cells.Item($row, 7).Value = $rec[4]
}
